$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the board-side markers in column D ("F.Cu" -> "Top", "B.Cu" -> "Bottom")
# for every data row (the pick-and-place side column).
$lastRow = 98
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $side = $cell.Value2
    if ($side -eq "F.Cu") {
        $cell.Value = "Top"
    } elseif ($side -eq "B.Cu") {
        $cell.Value = "Bottom"
    }
}

# Restore the view/selection state recorded in the saved workbook.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 69
$win.ScrollColumn = 1
$ws.Range("D99").Select()
